$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original 3 data rows (rows 2-4) shift down to rows 5-7; three new
# data rows are inserted above them as the new rows 2-4.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# All source values (invoice #, amounts, dates, currency codes, ...) are
# stored as plain text in this workbook, not as numbers/dates. Format the
# new block as text first so Excel doesn't auto-convert the numeric- and
# date-looking values (e.g. "574341", "2017-02-12") into real numbers.
$ws.Range("A2:G4").NumberFormat = "@"

# New row 2
$ws.Range("A2").Value = "574341"
$ws.Range("B2").Value = "IT Support"
$ws.Range("C2").Value = "2017-02-12"
$ws.Range("D2").Value = "201339"
$ws.Range("E2").Value = "40267.8"
$ws.Range("F2").Value = "241607"
$ws.Range("G2").Value = "USD"

# New row 3
$ws.Range("A3").Value = "544053"
$ws.Range("B3").Value = "Beverages and Catering"
$ws.Range("C3").Value = "2017-02-25"
$ws.Range("D3").Value = "270366"
$ws.Range("E3").Value = "54073.2"
$ws.Range("F3").Value = "324439"
$ws.Range("G3").Value = "USD"

# New row 4
$ws.Range("A4").Value = "830988"
$ws.Range("B4").Value = "Various paper supplies"
$ws.Range("C4").Value = "2017-02-24"
$ws.Range("D4").Value = "79397"
$ws.Range("E4").Value = "15879.4"
$ws.Range("F4").Value = "95276.4"
$ws.Range("G4").Value = "RON"
